$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
$np = $s1.NotesPage
$ph = $np.Shapes.Placeholders.Item(2)
$ph.TextFrame.TextRange.Text = "Speaker notes for slide one"
$nm = $p.NotesMaster
Write-Host "NM shapes:" $nm.Shapes.Count
